$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- New VLOOKUP demo data (rows 22-25, columns H-P) ---

# Row 22
$ws.Range("H22").Value = 1
$ws.Range("I22").Value = 10

# Row 23-25 helper values (column I)
$ws.Range("I23").Value = 20
$ws.Range("I24").Value = 30
$ws.Range("I25").Value = 40

# Names, in the order they were first entered (column J)
$ws.Range("J22").Value = "Paul"
$ws.Range("J23").Value = "John"
$ws.Range("J24").Value = "George"
$ws.Range("J25").Value = "Ringo"

# Instruments, entered bass, then lead guitar, then rythm guitar, then drums
$ws.Range("K22").Value = "bass"
$ws.Range("K24").Value = "lead guitar"
$ws.Range("K23").Value = "rythm guitar"
$ws.Range("K25").Value = "drums"

# Column L: L22 standalone formula, L23:L25 entered together as a shared formula
$ws.Range("L22").Formula = "=`$H`$22*20+I22"
$ws.Range("L23:L25").Formula = "=`$H`$22*20+I23"

# VLOOKUP examples
$ws.Range("N22").Formula = "=VLOOKUP(22,I22:L25,4)"
$ws.Range("O22").Formula = "=VLOOKUP(30,I22:L25,4,FALSE)"
$ws.Range("P22").Formula = "=VLOOKUP(31,I22:L25,4,FALSE)"

# --- View state: scroll down and select K28 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 8
$win.ScrollColumn = 1
$ws.Range("K28").Select()
